# UC007 - Listar Autorizações de Pagamento Pendentes
# The three "unique second step" cells (B/D on rows 20, 28, 36) rotate their
# content among the TC2, TC3 and TC4 blocks:
#   TC2 (row 20) gets what used to be TC4's step (atribuir/desatribuir)
#   TC3 (row 28) gets what used to be TC2's step (filtrar)
#   TC4 (row 36) gets what used to be TC3's step (realizar autorização)
# The "TCx" id labels (B15, B23, B31, B39) and every other cell stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$atribuir_step = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$atribuir_result = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

$filtrar_step = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$filtrar_result = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

$realizar_step = "Chefe Clica para realizar a autorização de pagamento."
$realizar_result = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"

# TC2 block (rows 15-20): second step becomes the "atribuir/desatribuir" content
$ws.Range("B20").Value = $atribuir_step
$ws.Range("D20").Value = $atribuir_result

# TC3 block (rows 23-28): second step becomes the "filtrar" content
$ws.Range("B28").Value = $filtrar_step
$ws.Range("D28").Value = $filtrar_result

# TC4 block (rows 31-36): second step becomes the "realizar autorização" content
$ws.Range("B36").Value = $realizar_step
$ws.Range("D36").Value = $realizar_result
